# Updated cryptos list (latest Price / Volume(1h) snapshot).
# The diff also re-orders two pairs of rows (39<->40 and 43<->44 swap places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A single apostrophe, built as a variable (not via [char] arithmetic, which
# this host coerces numerically when concatenated with a numeric-looking
# string).
$quote = "'"

# Helper: write a cell value. $AsText = $true means the value is numeric-
# looking (e.g. "313.99") but must stay a text cell, matching the source
# data (the Price column is stored as text in this workbook). Prefixing
# with a leading apostrophe forces that, the same way typing a quote-
# prefixed value into Excel by hand does.
function Set-CellValue {
    param($Ws, $Row, $Col, $Value, $AsText, $Quote)
    $cell = $Ws.Cells.Item($Row, $Col)
    if ($AsText -eq $true) {
        $cell.Value = $Quote + $Value
    } else {
        $cell.Value = $Value
    }
}

$colIndex = @{ B = 2; C = 3; D = 4; E = 5 }


# Row 2
Set-CellValue $ws 2 $colIndex.D '28.133.64' $false $quote
Set-CellValue $ws 2 $colIndex.E '  -1.30%  ' $false $quote

# Row 3
Set-CellValue $ws 3 $colIndex.D '1.793.25' $false $quote
Set-CellValue $ws 3 $colIndex.E '  -1.50%  ' $false $quote

# Row 4
Set-CellValue $ws 4 $colIndex.E '  +0.15%  ' $false $quote

# Row 5
Set-CellValue $ws 5 $colIndex.D '313.99' $true $quote
Set-CellValue $ws 5 $colIndex.E '  -0.39%  ' $false $quote

# Row 6
Set-CellValue $ws 6 $colIndex.E '  +0.07%  ' $false $quote

# Row 7
Set-CellValue $ws 7 $colIndex.D '0.5204' $true $quote
Set-CellValue $ws 7 $colIndex.E '  +2.18%  ' $false $quote

# Row 8
Set-CellValue $ws 8 $colIndex.E '  -3.46%  ' $false $quote

# Row 9
Set-CellValue $ws 9 $colIndex.D '0.07991' $true $quote
Set-CellValue $ws 9 $colIndex.E '  -2.18%  ' $false $quote

# Row 10
Set-CellValue $ws 10 $colIndex.E '  -0.57%  ' $false $quote

# Row 11
Set-CellValue $ws 11 $colIndex.E '  -1.32%  ' $false $quote

# Row 12
Set-CellValue $ws 12 $colIndex.D '6.289' $true $quote
Set-CellValue $ws 12 $colIndex.E '  -0.81%  ' $false $quote

# Row 13
Set-CellValue $ws 13 $colIndex.E '  +0.15%  ' $false $quote

# Row 14
Set-CellValue $ws 14 $colIndex.D '20.50' $true $quote
Set-CellValue $ws 14 $colIndex.E '  -2.88%  ' $false $quote

# Row 15
Set-CellValue $ws 15 $colIndex.D '7.285' $true $quote
Set-CellValue $ws 15 $colIndex.E '  -3.19%  ' $false $quote

# Row 16
Set-CellValue $ws 16 $colIndex.D '1.792.30' $false $quote
Set-CellValue $ws 16 $colIndex.E '  -1.26%  ' $false $quote

# Row 17
Set-CellValue $ws 17 $colIndex.D '91.72' $true $quote
Set-CellValue $ws 17 $colIndex.E '  -0.71%  ' $false $quote

# Row 18
Set-CellValue $ws 18 $colIndex.D '0.00001091' $true $quote
Set-CellValue $ws 18 $colIndex.E '  -3.93%  ' $false $quote

# Row 19
Set-CellValue $ws 19 $colIndex.D '0.06569' $true $quote
Set-CellValue $ws 19 $colIndex.E '  -1.33%  ' $false $quote

# Row 20
Set-CellValue $ws 20 $colIndex.E '  +0.16%  ' $false $quote

# Row 21
Set-CellValue $ws 21 $colIndex.D '17.32' $true $quote
Set-CellValue $ws 21 $colIndex.E '  -2.89%  ' $false $quote

# Row 22
Set-CellValue $ws 22 $colIndex.D '5.952' $true $quote
Set-CellValue $ws 22 $colIndex.E '  -2.40%  ' $false $quote

# Row 23
Set-CellValue $ws 23 $colIndex.D '28.175.35' $false $quote
Set-CellValue $ws 23 $colIndex.E '  -1.24%  ' $false $quote

# Row 24
Set-CellValue $ws 24 $colIndex.D '11.15' $true $quote
Set-CellValue $ws 24 $colIndex.E '  -2.43%  ' $false $quote

# Row 25
Set-CellValue $ws 25 $colIndex.E '  +0.26%  ' $false $quote

# Row 26
Set-CellValue $ws 26 $colIndex.D '160.63' $true $quote
Set-CellValue $ws 26 $colIndex.E '  +3.10%  ' $false $quote

# Row 27
Set-CellValue $ws 27 $colIndex.D '20.42' $true $quote
Set-CellValue $ws 27 $colIndex.E '  -4.26%  ' $false $quote

# Row 28
Set-CellValue $ws 28 $colIndex.D '1.994.64' $false $quote
Set-CellValue $ws 28 $colIndex.E '  -1.51%  ' $false $quote

# Row 29
Set-CellValue $ws 29 $colIndex.D '2.335' $true $quote
Set-CellValue $ws 29 $colIndex.E '  -2.77%  ' $false $quote

# Row 30
Set-CellValue $ws 30 $colIndex.E '  -2.41%  ' $false $quote

# Row 31
Set-CellValue $ws 31 $colIndex.D '0.1076' $true $quote
Set-CellValue $ws 31 $colIndex.E '  -1.93%  ' $false $quote

# Row 32
Set-CellValue $ws 32 $colIndex.D '1.051' $true $quote
Set-CellValue $ws 32 $colIndex.E '  -5.61%  ' $false $quote

# Row 33
Set-CellValue $ws 33 $colIndex.E '  +0.39%  ' $false $quote

# Row 34
Set-CellValue $ws 34 $colIndex.D '5.545' $true $quote
Set-CellValue $ws 34 $colIndex.E '  -3.93%  ' $false $quote

# Row 35
Set-CellValue $ws 35 $colIndex.D '0.07233' $true $quote
Set-CellValue $ws 35 $colIndex.E '  +2.39%  ' $false $quote

# Row 36
Set-CellValue $ws 36 $colIndex.D '12.03' $true $quote
Set-CellValue $ws 36 $colIndex.E '  +6.57%  ' $false $quote

# Row 37
Set-CellValue $ws 37 $colIndex.D '0.02305' $true $quote
Set-CellValue $ws 37 $colIndex.E '  -1.87%  ' $false $quote

# Row 38
Set-CellValue $ws 38 $colIndex.D '0.2143' $true $quote
Set-CellValue $ws 38 $colIndex.E '  -3.69%  ' $false $quote

# Row 39
Set-CellValue $ws 39 $colIndex.B 'FraxShare' $false $quote
Set-CellValue $ws 39 $colIndex.C 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' $false $quote
Set-CellValue $ws 39 $colIndex.D '8.671' $true $quote
Set-CellValue $ws 39 $colIndex.E '  -1.81%  ' $false $quote

# Row 40
Set-CellValue $ws 40 $colIndex.B 'InternetComputer(DFINITY)' $false $quote
Set-CellValue $ws 40 $colIndex.C 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' $false $quote
Set-CellValue $ws 40 $colIndex.D '5.059' $true $quote
Set-CellValue $ws 40 $colIndex.E '  -3.40%  ' $false $quote

# Row 41
Set-CellValue $ws 41 $colIndex.D '0.6159' $true $quote
Set-CellValue $ws 41 $colIndex.E '  -2.45%  ' $false $quote

# Row 42
Set-CellValue $ws 42 $colIndex.D '1.160' $true $quote
Set-CellValue $ws 42 $colIndex.E '  -1.66%  ' $false $quote

# Row 43
Set-CellValue $ws 43 $colIndex.B 'WEMIXTOKEN' $false $quote
Set-CellValue $ws 43 $colIndex.C 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' $false $quote
Set-CellValue $ws 43 $colIndex.D '1.337' $true $quote
Set-CellValue $ws 43 $colIndex.E '  -4.57%  ' $false $quote

# Row 44
Set-CellValue $ws 44 $colIndex.B 'EnergySwap' $false $quote
Set-CellValue $ws 44 $colIndex.C 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' $false $quote
Set-CellValue $ws 44 $colIndex.D '13.23' $true $quote
Set-CellValue $ws 44 $colIndex.E '  -1.65%  ' $false $quote

# Row 45
Set-CellValue $ws 45 $colIndex.D '3.760' $true $quote
Set-CellValue $ws 45 $colIndex.E '  +0.62%  ' $false $quote

# Row 46
Set-CellValue $ws 46 $colIndex.D '0.5954' $true $quote
Set-CellValue $ws 46 $colIndex.E '  +0.53%  ' $false $quote

# Row 47
Set-CellValue $ws 47 $colIndex.D '127.92' $true $quote
Set-CellValue $ws 47 $colIndex.E '  +2.31%  ' $false $quote

# Row 48
Set-CellValue $ws 48 $colIndex.D '1.219' $true $quote
Set-CellValue $ws 48 $colIndex.E '  +3.12%  ' $false $quote

# Row 49
Set-CellValue $ws 49 $colIndex.E '  -3.43%  ' $false $quote

# Row 50
Set-CellValue $ws 50 $colIndex.D '0.06744' $true $quote
Set-CellValue $ws 50 $colIndex.E '  -2.17%  ' $false $quote

# Row 51
Set-CellValue $ws 51 $colIndex.D '72.83' $true $quote
Set-CellValue $ws 51 $colIndex.E '  -1.99%  ' $false $quote
